$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 87) to the bottom of the profit log, matching
# the run performed on 2026-02-19.
$row = 87

# Column A holds the date as literal text (e.g. "02/19/2026"), just like
# every other row in the sheet. Pre-format the cell as Text so Excel does
# not auto-convert the date-looking string into a real date serial, then
# drop the formatting override so the cell is left with no explicit style
# (matching the rest of the data rows).
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "02/19/2026"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = 8992.99
$ws.Range("C$row").Value = 0.2474013351478005
$ws.Range("D$row").Value = 0.7525986648521995
$ws.Range("E$row").Value = -347.51
$ws.Range("F$row").Value = -39.27
$ws.Range("G$row").Value = -24173.65
$ws.Range("H$row").Value = -78.13
$ws.Range("I$row").Value = -1177.97
$ws.Range("J$row").Value = -34.62
$ws.Range("K$row").Value = -25351.62
$ws.Range("L$row").Value = -73.81999999999999
